$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to numeric-looking Price/Volume cells so Excel
# keeps them as literal text (matching the source inline strings) instead
# of auto-converting to numbers/percentages.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '309.79'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-3.21%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '50.46'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '3.43%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.171'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-1.43%'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-3.74%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.498'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-2.22%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.353'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '11.88%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.564'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-4.69%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1217'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-5.73%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1978'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.15%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09543'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '1.47%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.04737'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2.99%'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.66%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001279'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-4.04%'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-1.53%'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '2,012.39%'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.33%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.12%'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.05%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.999'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.20%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1358'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-1.79%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04158'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.27%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001272'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-2.74%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.003952'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-6.98%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001349'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.30%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02600'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '-4.15%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06043'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '6.57%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01127'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '78.54%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007884'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.46%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1424'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-1.26%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.008381'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '8.48%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007669'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-5.34%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3387'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '6.05%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00007336'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '6.00%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000749'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.25%'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-7.22%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002618'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-34.68%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002098'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.25%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0001999'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.25%'
